$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.3620791435241699
$ws.Range("B1").Value = 0.9812839031219482
$ws.Range("C1").Value = 4.841172695159912
$ws.Range("D1").Value = 1.783223390579224
$ws.Range("E1").Value = 1.010960578918457
